# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from serial date 45545 (2024-09-10) to 45546 (2024-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45545) {
        $cell.Value2 = 45546
    }
}
